# Daily attendance processing - reorder "Recorded By" author lists in column G
# so that "System" (and its duplicate "system") sort before the email address.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 7).End(-4162).Row  # xlUp = -4162

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G = 7 = "Recorded By"
    $val = $cell.Value2

    if ($val -eq "dnasr281@gmail.com, System") {
        $cell.Value2 = "System, dnasr281@gmail.com"
    }
    elseif ($val -eq "backup@backdoor.com, system, System") {
        $cell.Value2 = "backup@backdoor.com, System, system"
    }
}
